$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 - EOPD: claim "HiTechnicEOPD" (new shared string -> index 121)
$ws.Range("B20").Value = "HiTechnicEOPD"

# Row 17 - DPressure250: claim "DexterPressureSensor250" (new shared string -> index 122)
$ws.Range("B17").Value = "DexterPressureSensor250"

# Row 18 - DPressure500: claim "DexterPressureSensor500" (new shared string -> index 123)
$ws.Range("B18").Value = "DexterPressureSensor500"

# Mode column for rows 17/18 - "Pressure" (new shared string -> index 124)
$ws.Range("F17").Value = "Pressure"
$ws.Range("F18").Value = "Pressure"

# Mode column for row 20 - "Distance" (new shared string -> index 125)
$ws.Range("F20").Value = "Distance"

# Fill in the rest of the claim details (these already exist in the shared string table)
$ws.Range("C17").Value = "Y"
$ws.Range("D17").Value = "Lawrie"
$ws.Range("E17").Value = "N"
$ws.Range("G17").Value = "SampleProvider"

$ws.Range("C18").Value = "Y"
$ws.Range("D18").Value = "Lawrie"
$ws.Range("E18").Value = "N"
$ws.Range("G18").Value = "SampleProvider"

$ws.Range("C20").Value = "Y"
$ws.Range("D20").Value = "Lawrie"
$ws.Range("E20").Value = "N"
$ws.Range("G20").Value = "SampleProvider"

# Move the active selection on the frozen bottom-right pane to G22
$ws.Range("G22").Select()
